# Apply the "Fixses issues on instructor dashboard" edit:
#  - A8 text: FAILURE -> FAILURE#$%
#  - A11 gets a hyperlink (display text "!@#$%^") while the cell's own text
#    becomes "!#$%^&*("
#  - New row A12 with literal number 1234556
#  - Selection moves to D15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A8: FAILURE -> FAILURE#$%
$ws.Range("A8").Value = "FAILURE#$%"

# A11: add a hyperlink (its "TextToDisplay" becomes the cell text at first),
# then overwrite the cell text afterwards so the hyperlink's display
# attribute ("!@#$%^") and the cell's own text ("!#$%^&*() differ, matching
# the target workbook.
$ws.Hyperlinks.Add($ws.Range("A11"), "http://example.com", "", "", "!@#$%^")
$ws.Range("A11").Value = "!#$%^&*("

# New row 12: plain numeric literal
$ws.Range("A12").Value = 1234556

# Move the active selection to D15
$ws.Range("D15").Select()
